$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("展览")
$ws.Range("F4").Value = 1867
$ws.Range("F5").Value = 3253
$ws.Range("F7").Value = 4743
$ws.Range("F9").Value = 272
$ws.Range("F10").Value = 159
$ws.Range("F11").Value = 616
$ws.Range("F13").Value = 17
$ws.Range("F19").Value = 139
$ws.Range("F21").Value = 4700
$ws.Range("F22").Value = 13
$ws.Range("F23").Value = 32
$ws.Range("F25").Value = 5846
$ws.Range("F27").Value = 1184
$ws.Range("F29").Value = 654
$ws.Range("F30").Value = 4410
$ws.Range("F32").Value = 76
$ws.Range("F33").Value = 119
$ws.Range("F34").Value = 818
$ws.Range("F35").Value = 61
$ws.Range("F36").Value = 745
$ws.Range("C37").Value = "北京·第16届IJOY漫展XCGF游戏节"
$ws.Range("E37").Value = "2024.05.01 09:00-05.04 17:00"
$ws.Range("F37").Value = 767
$ws.Range("I37").Value = "//i2.hdslb.com/bfs/openplatform/202402/H86O2Jvv1707017473134.jpeg"

$ws = $wb.Worksheets.Item("演出")
$ws.Range("F3").Value = 33

$ws = $wb.Worksheets.Item("本地生活")
$ws.Range("F3").Value = 1079
$ws.Range("F4").Value = 34

$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F4").Value = 1079
$ws.Range("F5").Value = 34
$ws.Range("F7").Value = 1867
$ws.Range("F9").Value = 3253
$ws.Range("F11").Value = 4743
$ws.Range("F13").Value = 272
$ws.Range("F14").Value = 159
$ws.Range("F15").Value = 616
$ws.Range("F17").Value = 17
$ws.Range("F22").Value = 33
$ws.Range("F24").Value = 139
$ws.Range("F26").Value = 4700
$ws.Range("F27").Value = 13
$ws.Range("F28").Value = 32
$ws.Range("F30").Value = 5846
$ws.Range("F32").Value = 1184
$ws.Range("F34").Value = 654
$ws.Range("F35").Value = 4410
$ws.Range("F38").Value = 76
$ws.Range("F39").Value = 119
$ws.Range("F40").Value = 818
$ws.Range("F41").Value = 61
$ws.Range("F42").Value = 745
$ws.Range("C43").Value = "北京·第16届IJOY漫展XCGF游戏节"
$ws.Range("E43").Value = "2024.05.01 09:00-05.04 17:00"
$ws.Range("F43").Value = 767
$ws.Range("I43").Value = "//i2.hdslb.com/bfs/openplatform/202402/H86O2Jvv1707017473134.jpeg"
